$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("74:74").Insert()

$ws.Cells.Item(74, 1).Value = 4
$ws.Cells.Item(74, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(74, 3).Value = "Los Lagos"
$ws.Cells.Item(74, 4).Value = 44484
$ws.Cells.Item(74, 5).Value = 10
$ws.Cells.Item(74, 6).Value = 100112021
$ws.Cells.Item(74, 7).Value = "Ají"
$ws.Cells.Item(74, 8).Value = "Inferno"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 100
$ws.Cells.Item(74, 11).Value = 50000
$ws.Cells.Item(74, 12).Value = 50000
$ws.Cells.Item(74, 13).Value = 50000
$ws.Cells.Item(74, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(74, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(74, 16).Value = 4167
$ws.Cells.Item(74, 17).Value = 12
$ws.Cells.Item(74, 18).Value = "Hortaliza"
